# Apply updates to the "Artfynd" worksheet.
# Rows 15 and 16 swap their species-record data (columns A, D, E, F, G, H, Q, R),
# and the "Taxonsorteringsordning" value in column B is incremented by 1 for every
# affected row (8-22), landing on whichever row each record ends up in after the swap.
# Rows 19 and 20 undergo the same kind of swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple +1 increments to column B (Taxonsorteringsordning) ---
$ws.Range("B8").Value  = 79500
$ws.Range("B9").Value  = 91829
$ws.Range("B10").Value = 79715
$ws.Range("B11").Value = 83222
$ws.Range("B12").Value = 79715
$ws.Range("B13").Value = 83224
$ws.Range("B14").Value = 79715
$ws.Range("B17").Value = 80349
$ws.Range("B18").Value = 79715
$ws.Range("B21").Value = 83222
$ws.Range("B22").Value = 83222

# --- Row 15 becomes the former row 16 record (B incremented by 1) ---
$ws.Range("A15").Value = 130930230
$ws.Range("B15").Value = 78256
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 228579
$ws.Range("F15").Value = "Liten svartspik"
$ws.Range("G15").Value = "Chaenothecopsis nana"
$ws.Range("H15").Value = "Tibell"
$ws.Range("Q15").Value = 448404
$ws.Range("R15").Value = 7037411

# --- Row 16 becomes the former row 15 record (B incremented by 1) ---
$ws.Range("A16").Value = 130930219
$ws.Range("B16").Value = 92531
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 3298
$ws.Range("F16").Value = "Trådticka"
$ws.Range("G16").Value = "Climacocystis borealis"
$ws.Range("H16").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q16").Value = 448355
$ws.Range("R16").Value = 7037273

# --- Row 19 becomes the former row 20 record (B incremented by 1) ---
$ws.Range("A19").Value = 130930223
$ws.Range("B19").Value = 79715
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 1797
$ws.Range("F19").Value = "Mjölig dropplav"
$ws.Range("G19").Value = "Cliostomum leprosum"
$ws.Range("H19").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("Q19").Value = 448337
$ws.Range("R19").Value = 7037328

# --- Row 20 becomes the former row 19 record (B incremented by 1) ---
$ws.Range("A20").Value = 130930231
$ws.Range("B20").Value = 83224
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6440
$ws.Range("F20").Value = "Vitgrynig nållav"
$ws.Range("G20").Value = "Chaenotheca subroscida"
$ws.Range("H20").Value = "(Eitner) Zahlbr."
$ws.Range("Q20").Value = 448412
$ws.Range("R20").Value = 7037419
